# Removed Test Case Inter-Dependency
$wb = $excel.ActiveWorkbook

$wsInput  = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

$newProductName = "2475-RBI-EI-DB-SAR-REC-NOCOM-RNI-CTPD-DL-MD-TR-2-DATE-VAR-INST-PERIODIC-1st"
$newShortName   = "247e"

# Update the product name on both sheets (B1)
$wsInput.Range("B1").Value = $newProductName
$wsOutput.Range("B1").Value = $newProductName

# Update the short name on the input sheet (B2), now stored as text
$wsInput.Range("B2").Value = $newShortName

# Selection on input sheet moves back to B1
$wsInput.Range("B1").Select()

# Make the output sheet the active/selected sheet
$wsOutput.Activate()
$wsOutput.Range("B1").Select()
